{"js": "const replacements = [\n  [\"433\u00f79=48, 1\", \"847\u00f77=121, 0\"],\n  [\"713\u00f75=142, 3\", \"856\u00f75=171, 1\"],\n  [\"438\u00f73=146, 0\", \"683\u00f75=136, 3\"],\n  [\"448\u00f77=64, 0\", \"664\u00f77=94, 6\"],\n  [\"672\u00f76=112, 0\", \"175\u00f75=35, 0\"],\n  [\"434\u00f72=217, 0\", \"795\u00f72=397, 1\"],\n  [\"602\u00f74=150, 2\", \"591\u00f77=84, 3\"],\n  [\"384\u00f74=96, 0\", \"660\u00f79=73, 3\"],\n  [\"914\u00f77=130, 4\", \"971\u00f79=107, 8\"],\n  [\"494\u00f74=123, 2\", \"656\u00f74=164, 0\"],\n  [\"110\u00f76=18, 2\", \"943\u00f76=157, 1\"],\n  [\"327\u00f76=54, 3\", \"530\u00f72=265, 0\"],\n  [\"681\u00f78=85, 1\", \"132\u00f75=26, 2\"],\n  [\"234\u00f76=39, 0\", \"253\u00f76=42, 1\"],\n  [\"355\u00f79=39, 4\", \"948\u00f77=135, 3\"],\n  [\"179\u00f79=19, 8\", \"762\u00f72=381, 0\"],\n  [\"419\u00f73=139, 2\", \"456\u00f74=114, 0\"],\n  [\"939\u00f73=313, 0\", \"497\u00f75=99, 2\"],\n  [\"950\u00f77=135, 5\", \"797\u00f75=159, 2\"],\n  [\"810\u00f72=405, 0\", \"201\u00f73=67, 0\"],\n  [\"525\u00f72=262, 1\", \"540\u00f75=108, 0\"],\n  [\"906\u00f77=129, 3\", \"666\u00f78=83, 2\"],\n  [\"185\u00f75=37, 0\", \"795\u00f79=88, 3\"],\n  [\"984\u00f79=109, 3\", \"256\u00f74=64, 0\"],\n  [\"405\u00f72=202, 1\", \"510\u00f76=85, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  ,@(\"433\u00f79=48, 1\", \"847\u00f77=121, 0\")\n  ,@(\"713\u00f75=142, 3\", \"856\u00f75=171, 1\")\n  ,@(\"438\u00f73=146, 0\", \"683\u00f75=136, 3\")\n  ,@(\"448\u00f77=64, 0\", \"664\u00f77=94, 6\")\n  ,@(\"672\u00f76=112, 0\", \"175\u00f75=35, 0\")\n  ,@(\"434\u00f72=217, 0\", \"795\u00f72=397, 1\")\n  ,@(\"602\u00f74=150, 2\", \"591\u00f77=84, 3\")\n  ,@(\"384\u00f74=96, 0\", \"660\u00f79=73, 3\")\n  ,@(\"914\u00f77=130, 4\", \"971\u00f79=107, 8\")\n  ,@(\"494\u00f74=123, 2\", \"656\u00f74=164, 0\")\n  ,@(\"110\u00f76=18, 2\", \"943\u00f76=157, 1\")\n  ,@(\"327\u00f76=54, 3\", \"530\u00f72=265, 0\")\n  ,@(\"681\u00f78=85, 1\", \"132\u00f75=26, 2\")\n  ,@(\"234\u00f76=39, 0\", \"253\u00f76=42, 1\")\n  ,@(\"355\u00f79=39, 4\", \"948\u00f77=135, 3\")\n  ,@(\"179\u00f79=19, 8\", \"762\u00f72=381, 0\")\n  ,@(\"419\u00f73=139, 2\", \"456\u00f74=114, 0\")\n  ,@(\"939\u00f73=313, 0\", \"497\u00f75=99, 2\")\n  ,@(\"950\u00f77=135, 5\", \"797\u00f75=159, 2\")\n  ,@(\"810\u00f72=405, 0\", \"201\u00f73=67, 0\")\n  ,@(\"525\u00f72=262, 1\", \"540\u00f75=108, 0\")\n  ,@(\"906\u00f77=129, 3\", \"666\u00f78=83, 2\")\n  ,@(\"185\u00f75=37, 0\", \"795\u00f79=88, 3\")\n  ,@(\"984\u00f79=109, 3\", \"256\u00f74=64, 0\")\n  ,@(\"405\u00f72=202, 1\", \"510\u00f76=85, 0\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}"}
